$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.989.16'

$ws.Range("D3").Value = '2.205.55'
$ws.Range("E3").Value = '  +1.93%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '228.54'
$ws.Range("E5").Value = '  -0.23%  '

$ws.Range("E6").Value = '  +0.95%  '

$ws.Range("D7").Value = '63.54'
$ws.Range("E7").Value = '  +0.48%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -0.26%  '

$ws.Range("D10").Value = '0.0858'
$ws.Range("E10").Value = '  -0.58%  '

$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").Value = '16.08'
$ws.Range("E12").Value = '  +0.13%  '

$ws.Range("D13").Value = '2.531.65'
$ws.Range("E13").Value = '  +2.03%  '

$ws.Range("D14").Value = '22.15'
$ws.Range("E14").Value = '  -0.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.820'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("E16").Value = '  +0.42%  '

$ws.Range("D17").Value = '2.214.94'
$ws.Range("E17").Value = '  +2.41%  '

$ws.Range("D18").Value = '39.934.82'
$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("D19").Value = '0.0₃0918'
$ws.Range("E19").Value = '  +7.51%  '

$ws.Range("D20").Value = '72.19'
$ws.Range("E20").Value = '  -0.24%  '

$ws.Range("D21").Value = '6.09'
$ws.Range("E21").Value = '  -0.99%  '

$ws.Range("D22").Value = '231.73'
$ws.Range("E22").Value = '  +1.40%  '

$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '2.37'
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  -0.60%  '

$ws.Range("D26").Value = '9.56'
$ws.Range("E26").Value = '  -1.07%  '

$ws.Range("D27").Value = '171.34'
$ws.Range("E27").Value = '  -0.40%  '

$ws.Range("E28").Value = '  +1.42%  '

$ws.Range("D29").Value = '1.47'
$ws.Range("E29").Value = '  +3.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.16%  '

$ws.Range("E31").Value = '  +5.37%  '

$ws.Range("E32").Value = '  +0.74%  '

$ws.Range("E33").Value = '  -2.20%  '

$ws.Range("D34").Value = '4.73'
$ws.Range("E34").Value = '  -1.53%  '

$ws.Range("D35").Value = '7.05'
$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("E36").Value = '  +0.21%  '

$ws.Range("D37").Value = '3.88'
$ws.Range("E37").Value = '  +9.00%  '

$ws.Range("D38").Value = '2.47'
$ws.Range("E38").Value = '  +1.28%  '

$ws.Range("B39").Value = 'FTXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D39").Value = '5.06'
$ws.Range("E39").Value = '  +19.27%  '

$ws.Range("B40").Value = 'BinanceUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("D41").Value = '103.62'
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").Value = '17.83'
$ws.Range("E43").Value = '  -2.22%  '

$ws.Range("E44").Value = '  +3.18%  '

$ws.Range("D45").Value = '1.518.15'
$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("E46").Value = '  +3.28%  '

$ws.Range("E47").Value = '  +0.20%  '

$ws.Range("E48").Value = '  -0.79%  '

$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("E50").Value = '  +32.64%  '

$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = '49.98'
$ws.Range("E51").Value = '  +7.29%  '
